$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3: Priority stays 2, Browser changes to ChromeWeb (new string)
$ws.Range("B3").Value = "ChromeWeb"

# Update row 4: Priority changes 2 -> 4, Browser changes ChromeMobile -> SafariTablet
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "SafariTablet"

# Add new row 5: Priority 3, Browser ChromeMobile
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "ChromeMobile"

# Update selection to match target state
$ws.Range("A4").Select()
